# Update "paises.xlsx" (Pais sheet): refresh COVID case counters and swap the
# Etiopia / Republica de Macedonia rows back into alphabetical-ish order, plus
# bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Shared-string reorder effect: rows 79/80 swap which country they show.
# Row 79 used to be "Republica de Macedonia" and row 80 used to be "Etiopia";
# after the edit row 79 holds (updated) Etiopia data and row 80 holds the
# (unchanged) Republica de Macedonia data.
$ws.Range("A79").Value = "Etiopia"
$ws.Range("A80").Value = "Republica de Macedonia"

# --- Updated "last refreshed" banner text.
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 19:32"

# --- Per-country numeric updates (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes).

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3861874
$ws.Range("C4").Value = 28603
$ws.Range("D4").Value = 1776565
$ws.Range("E4").Value = 1942255
$ws.Range("G4").Value = 177
$ws.Range("H4").Value = 143054

# Brasil (row 5)
$ws.Range("B5").Value = 2076635
$ws.Range("C5").Value = 1389
$ws.Range("E5").Value = 630989
$ws.Range("G5").Value = 54
$ws.Range("H5").Value = 78871

# India (row 6)
$ws.Range("B6").Value = 1116613
$ws.Range("C6").Value = 38749
$ws.Range("D6").Value = 700305
$ws.Range("E6").Value = 388806
$ws.Range("G6").Value = 674
$ws.Range("H6").Value = 27502

# Turquia (row 18)
$ws.Range("B18").Value = 219641
$ws.Range("C18").Value = 924
$ws.Range("D18").Value = 202010
$ws.Range("E18").Value = 12140
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 5491

# Israel (row 44)
$ws.Range("B44").Value = 50035
$ws.Range("C44").Value = 670
$ws.Range("D44").Value = 21589
$ws.Range("E44").Value = 28037
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 409

# Irlanda (row 58)
$ws.Range("B58").Value = 25760
$ws.Range("C58").Value = 10
$ws.Range("E58").Value = 643

# Etiopia (row 79, after the name swap above)
$ws.Range("B79").Value = 9503
$ws.Range("C79").Value = 356
$ws.Range("D79").Value = 2430
$ws.Range("E79").Value = 6906
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 167

# Republica de Macedonia (row 80, after the name swap above)
$ws.Range("B80").Value = 9153
$ws.Range("C80").Value = 127
$ws.Range("D80").Value = 4810
$ws.Range("E80").Value = 3921
$ws.Range("G80").Value = 8
$ws.Range("H80").Value = 422

# Maldivas (row 109)
$ws.Range("B109").Value = 2966
$ws.Range("C109").Value = 36
$ws.Range("D109").Value = 2362
$ws.Range("E109").Value = 589

# Sri Lanka (row 113)
$ws.Range("B113").Value = 2724
$ws.Range("C113").Value = 20
$ws.Range("E113").Value = 678

# Cuba (row 116)
$ws.Range("B116").Value = 2446
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 2308
$ws.Range("E116").Value = 51

# Sudan del Sur (row 117)
$ws.Range("B117").Value = 2200
$ws.Range("C117").Value = 9
$ws.Range("E117").Value = 982

# Libia (row 127)
$ws.Range("B127").Value = 1866
$ws.Range("C127").Value = 75
$ws.Range("D127").Value = 418
$ws.Range("E127").Value = 1400

# Republica de Chipre (row 145)
$ws.Range("B145").Value = 1038
$ws.Range("C145").Value = 1
$ws.Range("E145").Value = 174

# Republica del Chad (row 148)
$ws.Range("D148").Value = 801
$ws.Range("E148").Value = 13

# Birmania (row 165)
$ws.Range("D165").Value = 276
$ws.Range("E165").Value = 59

# Comoras (row 167)
$ws.Range("B167").Value = 334
$ws.Range("C167").Value = 6
$ws.Range("D167").Value = 313
$ws.Range("E167").Value = 14
